# Fix(products): fix the sheet name and path create need fix data in xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Setembro"

# Delete column A, shifting B:AF left to A:AE
$ws.Columns.Item(1).Delete()
